# [fix] convert psi value deg to rad
#
# Column F ("psi_filter") was being accumulated from column G ("r_filter")
# without converting the per-step turn rate from degrees to radians before
# summing it into the running heading. Replace the hard-coded F3:F41 values
# with live formulas that wrap the G-column term in RADIANS(), matching the
# accumulation already used for F2 -> F3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F3 is the first row of the accumulation chain (references the static F2 seed).
$ws.Range("F3").Formula = "=F2+RADIANS(G2)+(A3-A2)"

# F4:F41 repeat the same relative pattern down the column (Excel adjusts the
# row-relative references automatically, just like a fill-down/autofill).
$ws.Range("F4:F41").Formula = "=F3+RADIANS(G3)+(A4-A3)"
